$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (NBOX STAR 40 HOURS ... item) price corrected to Rs. 749
$ws.Range("B3").Value = "Rs. 749"

# Row 4: item swapped from "NBOX INVICTUS ..." to the new "NBOX STAR WIRELESS NECKBAND ..." listing,
# with its price also set to Rs. 749
$ws.Range("A4").Value = "NBOX STAR WIRELESS NECKBAND WITH DOLBY EFFECT BASS SOUND IPX5 WITH MASSIVE MUSIC PLAYBACK WITH 1 YEAR WARRANTY BLUETOOTH HEADPHONE,BLUETOOTH EARPHONE,BLUETOOTH NECKBAND"
$ws.Range("B4").Value = "Rs. 749"
